$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Select whole rows 1:3 on ApplicationCheck (matches the stored selection state)
$ws1.Range("A1:XFD3").Select()

# Add new columns/data to ApplicationCheck sheet.
# NOTE: value-assignment order matters for the shared-strings table order,
# so the data rows are populated before the header row (matches source order).
$ws1.Range("D2").Value = "John"
$ws1.Range("E2").Value = "Doe"
$ws1.Range("F2").Value = 12345
$ws1.Range("G2").Value = "Sauce Labs Backpack"

$ws1.Range("D3").Value = "Jane"
$ws1.Range("E3").Value = "Smith"
$ws1.Range("F3").Value = 54321
$ws1.Range("G3").Value = "Sauce Labs Fleece Jacket"

$ws1.Range("D1").Value = "firstName"
$ws1.Range("E1").Value = "lastName"
$ws1.Range("F1").Value = "postalCode"
$ws1.Range("G1").Value = "productName"

$ws1.Columns.Item(7).ColumnWidth = 22.14
$ws1.PageSetup.Orientation = 1

# Add new sheet "SauceDemoShop" right after ApplicationCheck
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "SauceDemoShop"

$ws2.Range("A1").Value = "TestCase"
$ws2.Range("B1").Value = "username"
$ws2.Range("C1").Value = "password"
$ws2.Range("D1").Value = "firstName"
$ws2.Range("E1").Value = "lastName"
$ws2.Range("F1").Value = "postalCode"
$ws2.Range("G1").Value = "productName"

$ws2.Range("A2").Value = "Smoke"
$ws2.Range("B2").Value = "standard_user"
$ws2.Range("C2").Value = "secret_sauce"
$ws2.Range("D2").Value = "John"
$ws2.Range("E2").Value = "Doe"
$ws2.Range("F2").Value = 12345
$ws2.Range("G2").Value = "Sauce Labs Backpack"

$ws2.Range("A3").Value = "Regression"
$ws2.Range("B3").Value = "standard_user"
$ws2.Range("C3").Value = "secret_sauce"
$ws2.Range("D3").Value = "Jane"
$ws2.Range("E3").Value = "Smith"
$ws2.Range("F3").Value = 54321
$ws2.Range("G3").Value = "Sauce Labs Fleece Jacket"

$ws2.Range("E10").Select()
$ws2.PageSetup.Orientation = 1
